# Scheduled-runner market-data refresh: rewrite the raw price/profit
# columns (H:N) for the leve rows whose quoted marketboard data moved.
# Values come straight from the target OOXML -- no formulas live in
# these columns, so this is a plain value push per sheet/cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 312549.12
$ws.Range("I11").Value = 312549.12
$ws.Range("K11").Value = 312549.12
$ws.Range("M11").Value = -312409.12
$ws.Range("H17").Value = 1489.8334
$ws.Range("J17").Value = 1517.674
$ws.Range("L17").Value = 4553.022
$ws.Range("N17").Value = -4889.022
$ws.Range("H19").Value = 3493.6667
$ws.Range("I19").Value = 1994
$ws.Range("J19").Value = 4243.5
$ws.Range("K19").Value = 1994
$ws.Range("L19").Value = 4243.5
$ws.Range("M19").Value = -1819
$ws.Range("N19").Value = -4593.5
$ws.Range("H40").Value = 6998.2856
$ws.Range("I40").Value = 5572.125
$ws.Range("J40").Value = 8899.833000000001
$ws.Range("K40").Value = 5572.125
$ws.Range("L40").Value = 8899.833000000001
$ws.Range("M40").Value = -5397.125
$ws.Range("N40").Value = -9249.833000000001
$ws.Range("H58").Value = 3626.6924
$ws.Range("H62").Value = 62501836
$ws.Range("I62").Value = 62501836
$ws.Range("K62").Value = 62501836
$ws.Range("M62").Value = -62501212
$ws.Range("H65").Value = 62501836
$ws.Range("I65").Value = 62501836
$ws.Range("K65").Value = 312509180
$ws.Range("M65").Value = -312506060
$ws.Range("H106").Value = 3575.2727
$ws.Range("I106").Value = 5709.5
$ws.Range("K106").Value = 5709.5
$ws.Range("M106").Value = -5078.5
$ws.Range("H112").Value = 2622.6428
$ws.Range("J112").Value = 2670.1482
$ws.Range("L112").Value = 8010.444600000001
$ws.Range("N112").Value = -10226.4446
$ws.Range("H125").Value = 4489.2856
$ws.Range("I125").Value = 3872.25
$ws.Range("K125").Value = 34850.25
$ws.Range("M125").Value = -32390.25
$ws.Range("H127").Value = 43235.332
$ws.Range("I127").Value = 43235.332
$ws.Range("K127").Value = 129705.996
$ws.Range("M127").Value = -124745.996
$ws.Range("H131").Value = 2738.111
$ws.Range("I131").Value = 2347.6155
$ws.Range("J131").Value = 3753.4
$ws.Range("K131").Value = 7042.8465
$ws.Range("L131").Value = 11260.2
$ws.Range("M131").Value = -2002.8465
$ws.Range("N131").Value = -21340.2
$ws.Range("H137").Value = 3617.7
$ws.Range("I137").Value = 1779.8182
$ws.Range("J137").Value = 5864
$ws.Range("K137").Value = 5339.4546
$ws.Range("L137").Value = 17592
$ws.Range("M137").Value = -2789.4546
$ws.Range("N137").Value = -22692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2835.4473
$ws.Range("I32").Value = 2330.9858
$ws.Range("K32").Value = 2330.9858
$ws.Range("M32").Value = -2043.9858
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("H45").Value = 3196.24
$ws.Range("I45").Value = 2895.55
$ws.Range("K45").Value = 2895.55
$ws.Range("M45").Value = -2518.55
$ws.Range("H74").Value = 1607.0212
$ws.Range("I74").Value = 1211.9269
$ws.Range("J74").Value = 4306.8335
$ws.Range("K74").Value = 1211.9269
$ws.Range("L74").Value = 4306.8335
$ws.Range("M74").Value = -337.9268999999999
$ws.Range("N74").Value = -6054.8335
$ws.Range("H76").Value = 35144
$ws.Range("J76").Value = 35144
$ws.Range("L76").Value = 35144
$ws.Range("N76").Value = -35820
$ws.Range("H77").Value = 1607.0212
$ws.Range("I77").Value = 1211.9269
$ws.Range("J77").Value = 4306.8335
$ws.Range("K77").Value = 6059.6345
$ws.Range("L77").Value = 21534.1675
$ws.Range("M77").Value = -1691.6345
$ws.Range("N77").Value = -30270.1675
$ws.Range("H79").Value = 35144
$ws.Range("J79").Value = 35144
$ws.Range("L79").Value = 35144
$ws.Range("N79").Value = -37484
$ws.Range("H110").Value = 418339.25
$ws.Range("I110").Value = 501617.1
$ws.Range("K110").Value = 501617.1
$ws.Range("M110").Value = -499572.1
$ws.Range("N34").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 289015.28
$ws.Range("I31").Value = 1113254.6
$ws.Range("J31").Value = 3701.6538
$ws.Range("K31").Value = 1113254.6
$ws.Range("L31").Value = 3701.6538
$ws.Range("M31").Value = -1112959.6
$ws.Range("N31").Value = -4291.6538
$ws.Range("H34").Value = 289015.28
$ws.Range("I34").Value = 1113254.6
$ws.Range("J34").Value = 3701.6538
$ws.Range("K34").Value = 1113254.6
$ws.Range("L34").Value = 3701.6538
$ws.Range("M34").Value = -1113052.6
$ws.Range("N34").Value = -4105.6538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3171.075
$ws.Range("I131").Value = 1084
$ws.Range("J131").Value = 3866.7666
$ws.Range("K131").Value = 3252
$ws.Range("L131").Value = 11600.2998
$ws.Range("M131").Value = 1788
$ws.Range("N131").Value = -21680.2998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1450.2307
$ws.Range("I102").Value = 1478.25
$ws.Range("K102").Value = 1478.25
$ws.Range("M102").Value = 143.75
$ws.Range("H122").Value = 15333
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("H126").Value = 3177.5334
$ws.Range("I126").Value = 2285.4285
$ws.Range("J126").Value = 3958.125
$ws.Range("K126").Value = 6856.2855
$ws.Range("L126").Value = 11874.375
$ws.Range("M126").Value = -4386.2855
$ws.Range("N126").Value = -16814.375
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M122").ClearContents() | Out-Null
$ws.Range("N138").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1550
$ws.Range("J17").Value = 1550
$ws.Range("L17").Value = 1550
$ws.Range("N17").Value = -1890
$ws.Range("H46").Value = 5195.8335
$ws.Range("I46").Value = 4835
$ws.Range("J46").Value = 7000
$ws.Range("K46").Value = 4835
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = -4647
$ws.Range("N46").Value = -7376
$ws.Range("H122").Value = 1261616.5
$ws.Range("I122").Value = 1014785.1
$ws.Range("K122").Value = 3044355.3
$ws.Range("M122").Value = -3041905.3
$ws.Range("H125").Value = 80985
$ws.Range("J125").Value = 80985
$ws.Range("L125").Value = 80985
$ws.Range("N125").Value = -90825
$ws.Range("H136").Value = 2880.375
$ws.Range("I136").Value = 1925.6207
$ws.Range("J136").Value = 4337.6313
$ws.Range("K136").Value = 5776.8621
$ws.Range("L136").Value = 13012.8939
$ws.Range("M136").Value = -3226.8621
$ws.Range("N136").Value = -18112.8939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14899.6
$ws.Range("J45").Value = 14666.333
$ws.Range("L45").Value = 14666.333
$ws.Range("N45").Value = -15648.333
$ws.Range("H49").Value = 10000
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9770
$ws.Range("H81").Value = 8152.2354
$ws.Range("I81").Value = 1059.5
$ws.Range("J81").Value = 18284.715
$ws.Range("K81").Value = 2119
$ws.Range("L81").Value = 36569.43
$ws.Range("M81").Value = -1058
$ws.Range("N81").Value = -38691.43
$ws.Range("H84").Value = 8152.2354
$ws.Range("I84").Value = 1059.5
$ws.Range("J84").Value = 18284.715
$ws.Range("K84").Value = 10595
$ws.Range("L84").Value = 182847.15
$ws.Range("M84").Value = -5291
$ws.Range("N84").Value = -193455.15
$ws.Range("H113").Value = 1463.24
$ws.Range("I113").Value = 1136.6316
$ws.Range("K113").Value = 3409.8948
$ws.Range("M113").Value = -1239.8948
$ws.Range("H122").Value = 29416368
$ws.Range("I122").Value = 47622710
$ws.Range("K122").Value = 142868130
$ws.Range("M122").Value = -142865680
